$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: capture all existing hyperlinks (anchor row/col, target address, display text) ---
$hlRows = @()
$hlCols = @()
$hlTargets = @()
$hlDisps = @()

foreach ($hl in $ws.Hyperlinks) {
    $hlRows += $hl.Range.Row
    $hlCols += $hl.Range.Column
    $hlTargets += $hl.Address
    $hlDisps += $hl.TextToDisplay
}

# --- Step 2: remove all hyperlinks so they don't interfere with the row delete ---
$ws.Hyperlinks.Delete()

# --- Step 3: delete the entire row 23 (shifts rows 24+ up by one) ---
$ws.Rows("23:23").Delete()

# --- Step 4: re-create the hyperlinks, shifting anchors that were below the deleted row.
#     Adding a hyperlink via COM auto-applies the built-in "Hyperlink" cell style, which
#     would corrupt the original cell formatting/style index, so we snapshot + restore the
#     individual font properties around the Add() call to keep the original style intact. ---
for ($i = 0; $i -lt $hlRows.Count; $i++) {
    $r = $hlRows[$i]
    if ($r -eq 23) {
        continue
    }
    if ($r -gt 23) {
        $r = $r - 1
    }
    $target = $ws.Cells.Item($r, $hlCols[$i])

    $fname = $target.Font.Name
    $fsize = $target.Font.Size
    $fcolor = $target.Font.Color
    $funderline = $target.Font.Underline
    $fbold = $target.Font.Bold
    $fitalic = $target.Font.Italic

    $ws.Hyperlinks.Add($target, $hlTargets[$i], "", "", $hlDisps[$i])

    $target.Font.Name = $fname
    $target.Font.Size = $fsize
    $target.Font.Color = $fcolor
    $target.Font.Underline = $funderline
    $target.Font.Bold = $fbold
    $target.Font.Italic = $fitalic
}

# --- Step 5: restore the recorded selection state ---
$ws.Range("A23").Select()
